# chore: adapt column header formatting to respective input file names
#
# 1. Rename the "_old" / "_new" column-header suffixes to the
#    format-version-specific "_FV2304" / "_FV2310" suffixes.
# 2. Turn the data range into a proper Excel Table ("Table1").
# 3. Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames ----------------------------------------------------
$oldSuffix = "_old"
$newSuffix = "_new"
$fv2304 = "_FV2304"
$fv2310 = "_FV2310"

$headerRange = $ws.Range("A1:U1")
for ($i = 1; $i -le $headerRange.Columns.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $text = [string]$cell.Value()
    if ($text.EndsWith($oldSuffix)) {
        $base = $text.Substring(0, $text.Length - $oldSuffix.Length)
        $cell.Value = $base + $fv2304
    } elseif ($text.EndsWith($newSuffix)) {
        $base = $text.Substring(0, $text.Length - $newSuffix.Length)
        $cell.Value = $base + $fv2310
    }
}

# --- 2. Convert the used range into an Excel Table ------------------------
$dataRange = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
